# "vylepseni v ip setting" - update the saved IP-address list and default
# interface setting.
#
# Sheet "ip_address_list": the placeholder row 1 ("projekt" / poznvv / ...)
# is removed, so the previous rows 2 and 3 shift up to become rows 1 and 2.
$wb = $excel.ActiveWorkbook

$wsAddr = $wb.Worksheets.Item("ip_address_list")
$wsAddr.Rows(1).Delete()

# Sheet "Settings": default interface index changes from 5 to 1.
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 1
